# "updated activity till excel form"
#
# The sheet tracks Jos Buttler's innings-by-innings activity (runs, balls,
# fours, sixes). A new innings row is inserted as row 4, pushing the
# previously-last two innings down by one row: what used to be row 4 becomes
# row 6, row 5's figures move up into row 4 with freshly reported numbers,
# and row 5 itself gets the newest innings' figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0) - the existing "Jos Buttler" entries in column A
# use it after the name, so keep it consistent for the new row.
$nbsp = [char]0x00A0

# Row 6 (new last row): carries the original row 4 figures (24 / 25 / 1 / 1).
$ws.Range("A6").Value = "Jos Buttler" + $nbsp
$ws.Range("B6").Value = "Rajasthan Royals"
$ws.Range("C6:F6").NumberFormat = "@"
$ws.Range("C6").Value = "24"
$ws.Range("D6").Value = "25"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "1"
# Drop the one-off text NumberFormat back to Normal so the cell keeps its
# text-stored-as-number value without leaving a stray style behind (matches
# the rest of the sheet, which has no explicit per-cell formatting).
$ws.Range("C6:F6").Style = "Normal"

# Row 4: now holds what used to be row 5's figures (35 / 22 / 4 / 1).
$ws.Range("C4:F4").NumberFormat = "@"
$ws.Range("C4").Value = "35"
$ws.Range("D4").Value = "22"
$ws.Range("E4").Value = "4"
$ws.Range("F4").Value = "1"
$ws.Range("C4:F4").Style = "Normal"

# Row 5: the newly recorded innings (70 / 48 / 7 / 2).
$ws.Range("C5:F5").NumberFormat = "@"
$ws.Range("C5").Value = "70"
$ws.Range("D5").Value = "48"
$ws.Range("E5").Value = "7"
$ws.Range("F5").Value = "2"
$ws.Range("C5:F5").Style = "Normal"
